$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows of Mac-Address (machine) master data, following the
# existing pattern: regcntr_id, machine_id, lang_code, is_active, cr_by, cr_dtimes
$ws.Range("A31").Value = 10001
$ws.Range("B31").Value = 10030
$ws.Range("C31").Value = "eng"
$ws.Range("D31").Value = $true
$ws.Range("E31").Value = "superadmin"
$ws.Range("F31").Value = "now()"

$ws.Range("A32").Value = 10001
$ws.Range("B32").Value = 10031
$ws.Range("C32").Value = "eng"
$ws.Range("D32").Value = $true
$ws.Range("E32").Value = "superadmin"
$ws.Range("F32").Value = "now()"

# Update the selection/view to match: scroll to top, select E31
$ws.Range("E31").Select()
